$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 7).Value = 66.47695399999999
$ws.Cells.Item(2, 8).Value = 199.430862
$ws.Cells.Item(2, 9).Value = 0.04311983106164722
$ws.Cells.Item(2, 10).Value = 0.04311983106164721
$ws.Cells.Item(2, 13).Value = 0.3331066666666667
$ws.Cells.Item(2, 14).Value = 0.99932
$ws.Cells.Item(2, 15).Value = 0.002125805913843485
$ws.Cells.Item(2, 16).Value = 0.002125805913843485
$ws.Cells.Item(2, 17).Value = 22.14391655709333
$ws.Cells.Item(2, 18).Value = 199.29524901384
$ws.Cells.Item(2, 19).Value = 0.00009166439187478166
$ws.Cells.Item(2, 20).Value = 0.00009166439187478164
$ws.Cells.Item(3, 7).Value = 66.47695399999999
$ws.Cells.Item(3, 8).Value = 199.430862
$ws.Cells.Item(3, 9).Value = 0.04311983106164722
$ws.Cells.Item(3, 10).Value = 0.04311983106164721
$ws.Cells.Item(3, 15).Value = 0.0008775937418887864
$ws.Cells.Item(3, 16).Value = 0.0008775937418887864
$ws.Cells.Item(3, 17).Value = 9.141644806263999
$ws.Cells.Item(3, 18).Value = 82.274803256376
$ws.Cells.Item(3, 19).Value = 0.0000378416938910033
$ws.Cells.Item(3, 20).Value = 0.0000378416938910033
$ws.Cells.Item(4, 7).Value = 66.47695399999999
$ws.Cells.Item(4, 8).Value = 199.430862
$ws.Cells.Item(4, 9).Value = 0.04311983106164722
$ws.Cells.Item(4, 10).Value = 0.04311983106164721
$ws.Cells.Item(4, 13).Value = 91.40156066666667
$ws.Cells.Item(4, 14).Value = 274.204682
$ws.Cells.Item(4, 15).Value = 0.5833025803538128
$ws.Cells.Item(4, 16).Value = 0.5833025803538128
$ws.Cells.Item(4, 17).Value = 6076.097343966208
$ws.Cells.Item(4, 18).Value = 54684.87609569588
$ws.Cells.Item(4, 19).Value = 0.02515190872267931
$ws.Cells.Item(4, 20).Value = 0.02515190872267931
$ws.Cells.Item(5, 7).Value = 66.47695399999999
$ws.Cells.Item(5, 8).Value = 199.430862
$ws.Cells.Item(5, 9).Value = 0.04311983106164722
$ws.Cells.Item(5, 10).Value = 0.04311983106164721
$ws.Cells.Item(5, 13).Value = 0.5759770000000001
$ws.Cells.Item(5, 14).Value = 1.727931
$ws.Cells.Item(5, 15).Value = 0.00367574544541637
$ws.Cells.Item(5, 16).Value = 0.00367574544541637
$ws.Cells.Item(5, 17).Value = 38.289196534058
$ws.Cells.Item(5, 18).Value = 344.602768806522
$ws.Cells.Item(5, 19).Value = 0.0001584975226319731
$ws.Cells.Item(5, 20).Value = 0.0001584975226319731
$ws.Cells.Item(6, 7).Value = 66.47695399999999
$ws.Cells.Item(6, 8).Value = 199.430862
$ws.Cells.Item(6, 9).Value = 0.04311983106164722
$ws.Cells.Item(6, 10).Value = 0.04311983106164721
$ws.Cells.Item(6, 13).Value = 64.24849033333334
$ws.Cells.Item(6, 14).Value = 192.745471
$ws.Cells.Item(6, 15).Value = 0.4100182745450386
$ws.Cells.Item(6, 16).Value = 0.4100182745450385
$ws.Cells.Item(6, 17).Value = 4271.043936458444
$ws.Cells.Item(6, 18).Value = 38439.395428126
$ws.Cells.Item(6, 19).Value = 0.01767991873057015
$ws.Cells.Item(6, 20).Value = 0.01767991873057015
$ws.Cells.Item(7, 9).Value = 0.8830494168872806
$ws.Cells.Item(7, 10).Value = 0.8830494168872804
$ws.Cells.Item(7, 13).Value = 0.3331066666666667
$ws.Cells.Item(7, 14).Value = 0.99932
$ws.Cells.Item(7, 15).Value = 0.002125805913843485
$ws.Cells.Item(7, 16).Value = 0.002125805913843485
$ws.Cells.Item(7, 17).Value = 453.48444374436
$ws.Cells.Item(7, 18).Value = 4081.35999369924
$ws.Cells.Item(7, 19).Value = 0.001877191672635022
$ws.Cells.Item(7, 20).Value = 0.001877191672635022
$ws.Cells.Item(8, 9).Value = 0.8830494168872806
$ws.Cells.Item(8, 10).Value = 0.8830494168872804
$ws.Cells.Item(8, 15).Value = 0.0008775937418887864
$ws.Cells.Item(8, 16).Value = 0.0008775937418887864
$ws.Cells.Item(8, 19).Value = 0.0007749586420388194
$ws.Cells.Item(8, 20).Value = 0.0007749586420388193
$ws.Cells.Item(9, 9).Value = 0.8830494168872806
$ws.Cells.Item(9, 10).Value = 0.8830494168872804
$ws.Cells.Item(9, 13).Value = 91.40156066666667
$ws.Cells.Item(9, 14).Value = 274.204682
$ws.Cells.Item(9, 15).Value = 0.5833025803538128
$ws.Cells.Item(9, 16).Value = 0.5833025803538128
$ws.Cells.Item(9, 17).Value = 124432.1715655337
$ws.Cells.Item(9, 18).Value = 1119889.544089803
$ws.Cells.Item(9, 19).Value = 0.5150850034502805
$ws.Cells.Item(9, 20).Value = 0.5150850034502804
$ws.Cells.Item(10, 9).Value = 0.8830494168872806
$ws.Cells.Item(10, 10).Value = 0.8830494168872804
$ws.Cells.Item(10, 13).Value = 0.5759770000000001
$ws.Cells.Item(10, 14).Value = 1.727931
$ws.Cells.Item(10, 15).Value = 0.00367574544541637
$ws.Cells.Item(10, 16).Value = 0.00367574544541637
$ws.Cells.Item(10, 17).Value = 784.123032025413
$ws.Cells.Item(10, 18).Value = 7057.107288228717
$ws.Cells.Item(10, 19).Value = 0.003245864872201003
$ws.Cells.Item(10, 20).Value = 0.003245864872201002
$ws.Cells.Item(11, 9).Value = 0.8830494168872806
$ws.Cells.Item(11, 10).Value = 0.8830494168872804
$ws.Cells.Item(11, 13).Value = 64.24849033333334
$ws.Cells.Item(11, 14).Value = 192.745471
$ws.Cells.Item(11, 15).Value = 0.4100182745450386
$ws.Cells.Item(11, 16).Value = 0.4100182745450385
$ws.Cells.Item(11, 17).Value = 87466.54995464883
$ws.Cells.Item(11, 18).Value = 787198.9495918395
$ws.Cells.Item(11, 19).Value = 0.3620663982501252
$ws.Cells.Item(11, 20).Value = 0.3620663982501252
$ws.Cells.Item(12, 7).Value = 44.831112
$ws.Cells.Item(12, 8).Value = 134.493336
$ws.Cells.Item(12, 9).Value = 0.02907940059566787
$ws.Cells.Item(12, 10).Value = 0.02907940059566786
$ws.Cells.Item(12, 13).Value = 0.3331066666666667
$ws.Cells.Item(12, 14).Value = 0.99932
$ws.Cells.Item(12, 15).Value = 0.002125805913843485
$ws.Cells.Item(12, 16).Value = 0.002125805913843485
$ws.Cells.Item(12, 17).Value = 14.93354228128
$ws.Cells.Item(12, 18).Value = 134.40188053152
$ws.Cells.Item(12, 19).Value = 0.00006181716175729452
$ws.Cells.Item(12, 20).Value = 0.00006181716175729451
$ws.Cells.Item(13, 7).Value = 44.831112
$ws.Cells.Item(13, 8).Value = 134.493336
$ws.Cells.Item(13, 9).Value = 0.02907940059566787
$ws.Cells.Item(13, 10).Value = 0.02907940059566786
$ws.Cells.Item(13, 15).Value = 0.0008775937418887864
$ws.Cells.Item(13, 16).Value = 0.0008775937418887864
$ws.Cells.Item(13, 17).Value = 6.164995197792
$ws.Cells.Item(13, 18).Value = 55.484956780128
$ws.Cells.Item(13, 19).Value = 0.00002551989998063517
$ws.Cells.Item(13, 20).Value = 0.00002551989998063516
$ws.Cells.Item(14, 7).Value = 44.831112
$ws.Cells.Item(14, 8).Value = 134.493336
$ws.Cells.Item(14, 9).Value = 0.02907940059566787
$ws.Cells.Item(14, 10).Value = 0.02907940059566786
$ws.Cells.Item(14, 13).Value = 91.40156066666667
$ws.Cells.Item(14, 14).Value = 274.204682
$ws.Cells.Item(14, 15).Value = 0.5833025803538128
$ws.Cells.Item(14, 16).Value = 0.5833025803538128
$ws.Cells.Item(14, 17).Value = 4097.633603222128
$ws.Cells.Item(14, 18).Value = 36878.70242899915
$ws.Cells.Item(14, 19).Value = 0.01696208940259527
$ws.Cells.Item(14, 20).Value = 0.01696208940259527
$ws.Cells.Item(15, 7).Value = 44.831112
$ws.Cells.Item(15, 8).Value = 134.493336
$ws.Cells.Item(15, 9).Value = 0.02907940059566787
$ws.Cells.Item(15, 10).Value = 0.02907940059566786
$ws.Cells.Item(15, 13).Value = 0.5759770000000001
$ws.Cells.Item(15, 14).Value = 1.727931
$ws.Cells.Item(15, 15).Value = 0.00367574544541637
$ws.Cells.Item(15, 16).Value = 0.00367574544541637
$ws.Cells.Item(15, 17).Value = 25.821689396424
$ws.Cells.Item(15, 18).Value = 232.395204567816
$ws.Cells.Item(15, 19).Value = 0.0001068884742949643
$ws.Cells.Item(15, 20).Value = 0.0001068884742949642
$ws.Cells.Item(16, 7).Value = 44.831112
$ws.Cells.Item(16, 8).Value = 134.493336
$ws.Cells.Item(16, 9).Value = 0.02907940059566787
$ws.Cells.Item(16, 10).Value = 0.02907940059566786
$ws.Cells.Item(16, 13).Value = 64.24849033333334
$ws.Cells.Item(16, 14).Value = 192.745471
$ws.Cells.Item(16, 15).Value = 0.4100182745450386
$ws.Cells.Item(16, 16).Value = 0.4100182745450385
$ws.Cells.Item(16, 17).Value = 2880.331265964584
$ws.Cells.Item(16, 18).Value = 25922.98139368126
$ws.Cells.Item(16, 19).Value = 0.01192308565703971
$ws.Cells.Item(16, 20).Value = 0.0119230856570397
$ws.Cells.Item(17, 7).Value = 52.83062100000001
$ws.Cells.Item(17, 8).Value = 158.491863
$ws.Cells.Item(17, 9).Value = 0.0342682285413064
$ws.Cells.Item(17, 10).Value = 0.03426822854130639
$ws.Cells.Item(17, 13).Value = 0.3331066666666667
$ws.Cells.Item(17, 14).Value = 0.99932
$ws.Cells.Item(17, 15).Value = 0.002125805913843485
$ws.Cells.Item(17, 16).Value = 0.002125805913843485
$ws.Cells.Item(17, 17).Value = 17.59823205924
$ws.Cells.Item(17, 18).Value = 158.38408853316
$ws.Cells.Item(17, 19).Value = 0.00007284760289004925
$ws.Cells.Item(17, 20).Value = 0.00007284760289004923
$ws.Cells.Item(18, 7).Value = 52.83062100000001
$ws.Cells.Item(18, 8).Value = 158.491863
$ws.Cells.Item(18, 9).Value = 0.0342682285413064
$ws.Cells.Item(18, 10).Value = 0.03426822854130639
$ws.Cells.Item(18, 15).Value = 0.0008775937418887864
$ws.Cells.Item(18, 16).Value = 0.0008775937418887864
$ws.Cells.Item(18, 17).Value = 7.265055677436001
$ws.Cells.Item(18, 18).Value = 65.38550109692402
$ws.Cells.Item(18, 19).Value = 0.00003007358291346519
$ws.Cells.Item(18, 20).Value = 0.00003007358291346519
$ws.Cells.Item(19, 7).Value = 52.83062100000001
$ws.Cells.Item(19, 8).Value = 158.491863
$ws.Cells.Item(19, 9).Value = 0.0342682285413064
$ws.Cells.Item(19, 10).Value = 0.03426822854130639
$ws.Cells.Item(19, 13).Value = 91.40156066666667
$ws.Cells.Item(19, 14).Value = 274.204682
$ws.Cells.Item(19, 15).Value = 0.5833025803538128
$ws.Cells.Item(19, 16).Value = 0.5833025803538128
$ws.Cells.Item(19, 17).Value = 4828.801210389175
$ws.Cells.Item(19, 18).Value = 43459.21089350257
$ws.Cells.Item(19, 19).Value = 0.0199887461322982
$ws.Cells.Item(19, 20).Value = 0.01998874613229819
$ws.Cells.Item(20, 7).Value = 52.83062100000001
$ws.Cells.Item(20, 8).Value = 158.491863
$ws.Cells.Item(20, 9).Value = 0.0342682285413064
$ws.Cells.Item(20, 10).Value = 0.03426822854130639
$ws.Cells.Item(20, 13).Value = 0.5759770000000001
$ws.Cells.Item(20, 14).Value = 1.727931
$ws.Cells.Item(20, 15).Value = 0.00367574544541637
$ws.Cells.Item(20, 16).Value = 0.00367574544541637
$ws.Cells.Item(20, 17).Value = 30.42922259171701
$ws.Cells.Item(20, 18).Value = 273.8630033254531
$ws.Cells.Item(20, 19).Value = 0.0001259612849831943
$ws.Cells.Item(20, 20).Value = 0.0001259612849831942
$ws.Cells.Item(21, 7).Value = 52.83062100000001
$ws.Cells.Item(21, 8).Value = 158.491863
$ws.Cells.Item(21, 9).Value = 0.0342682285413064
$ws.Cells.Item(21, 10).Value = 0.03426822854130639
$ws.Cells.Item(21, 13).Value = 64.24849033333334
$ws.Cells.Item(21, 14).Value = 192.745471
$ws.Cells.Item(21, 15).Value = 0.4100182745450386
$ws.Cells.Item(21, 16).Value = 0.4100182745450385
$ws.Cells.Item(21, 17).Value = 3394.287642622498
$ws.Cells.Item(21, 18).Value = 30548.58878360248
$ws.Cells.Item(21, 19).Value = 0.01405059993822149
$ws.Cells.Item(21, 20).Value = 0.01405059993822149
$ws.Cells.Item(22, 7).Value = 16.16161433333333
$ws.Cells.Item(22, 8).Value = 48.484843
$ws.Cells.Item(22, 9).Value = 0.01048312291409786
$ws.Cells.Item(22, 10).Value = 0.01048312291409786
$ws.Cells.Item(22, 13).Value = 0.3331066666666667
$ws.Cells.Item(22, 14).Value = 0.99932
$ws.Cells.Item(22, 15).Value = 0.002125805913843485
$ws.Cells.Item(22, 16).Value = 0.002125805913843485
$ws.Cells.Item(22, 17).Value = 5.383541478528889
$ws.Cells.Item(22, 18).Value = 48.45187330676
$ws.Cells.Item(22, 19).Value = 0.00002228508468633739
$ws.Cells.Item(22, 20).Value = 0.00002228508468633738
$ws.Cells.Item(23, 7).Value = 16.16161433333333
$ws.Cells.Item(23, 8).Value = 48.484843
$ws.Cells.Item(23, 9).Value = 0.01048312291409786
$ws.Cells.Item(23, 10).Value = 0.01048312291409786
$ws.Cells.Item(23, 15).Value = 0.0008775937418887864
$ws.Cells.Item(23, 16).Value = 0.0008775937418887864
$ws.Cells.Item(23, 17).Value = 2.222480556662667
$ws.Cells.Item(23, 18).Value = 20.002325009964
$ws.Cells.Item(23, 19).Value = 0.000009199923064863223
$ws.Cells.Item(23, 20).Value = 0.000009199923064863221
$ws.Cells.Item(24, 7).Value = 16.16161433333333
$ws.Cells.Item(24, 8).Value = 48.484843
$ws.Cells.Item(24, 9).Value = 0.01048312291409786
$ws.Cells.Item(24, 10).Value = 0.01048312291409786
$ws.Cells.Item(24, 13).Value = 91.40156066666667
$ws.Cells.Item(24, 14).Value = 274.204682
$ws.Cells.Item(24, 15).Value = 0.5833025803538128
$ws.Cells.Item(24, 16).Value = 0.5833025803538128
$ws.Cells.Item(24, 17).Value = 1477.196772959436
$ws.Cells.Item(24, 18).Value = 13294.77095663493
$ws.Cells.Item(24, 19).Value = 0.006114832645959465
$ws.Cells.Item(24, 20).Value = 0.006114832645959464
$ws.Cells.Item(25, 7).Value = 16.16161433333333
$ws.Cells.Item(25, 8).Value = 48.484843
$ws.Cells.Item(25, 9).Value = 0.01048312291409786
$ws.Cells.Item(25, 10).Value = 0.01048312291409786
$ws.Cells.Item(25, 13).Value = 0.5759770000000001
$ws.Cells.Item(25, 14).Value = 1.727931
$ws.Cells.Item(25, 15).Value = 0.00367574544541637
$ws.Cells.Item(25, 16).Value = 0.00367574544541637
$ws.Cells.Item(25, 17).Value = 9.308718138870335
$ws.Cells.Item(25, 18).Value = 83.778463249833
$ws.Cells.Item(25, 19).Value = 0.00003853329130523521
$ws.Cells.Item(25, 20).Value = 0.0000385332913052352
$ws.Cells.Item(26, 7).Value = 16.16161433333333
$ws.Cells.Item(26, 8).Value = 48.484843
$ws.Cells.Item(26, 9).Value = 0.01048312291409786
$ws.Cells.Item(26, 10).Value = 0.01048312291409786
$ws.Cells.Item(26, 13).Value = 64.24849033333334
$ws.Cells.Item(26, 14).Value = 192.745471
$ws.Cells.Item(26, 15).Value = 0.4100182745450386
$ws.Cells.Item(26, 16).Value = 0.4100182745450385
$ws.Cells.Item(26, 17).Value = 1038.359322266228
$ws.Cells.Item(26, 18).Value = 9345.233900396053
$ws.Cells.Item(26, 19).Value = 0.004298271969081963
$ws.Cells.Item(26, 20).Value = 0.004298271969081962
